# "added user magement section"
# Adds a new row to the status list (B9 = "good") and highlights the
# header cell B1 with a distinct font (JetBrains Mono, green, vertically
# centered) to call out the new user-management section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New status entry in column B.
$ws.Range("B9").Value = "good"

# Restyle the "Erayga_Asalka" header cell (B1) so it stands out.
$hdr = $ws.Range("B1")
$hdr.Font.Name = "JetBrains Mono"
$hdr.Font.Size = 9.8
$hdr.Font.Color = 1539334        # RGB(6,125,23) -> 0xFF067D17
$hdr.VerticalAlignment = -4108   # xlCenter

# Leave the selection on the header cell, matching the author's last view.
$hdr.Select()
